# Update "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I) on the
# Training Dashboard sheet to reflect progress as of 04-Nov-2025
# (previously 03-Nov-2025): each period-to-expire value decreases by 1
# day and the last-update date moves forward to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$rows = 3, 4, 5, 6, 7

# A helper cell (outside the used range) is used to stage the new date
# text with a "Text" number format so that Excel does not auto-convert
# the literal "04-Nov-2025" string into a date serial number when it is
# pasted into column I (the source workbook stores these dates as plain
# text, not as real date values).
$helper = $ws.Cells.Item(20, 1)
$helper.NumberFormat = "@"
$helper.Value = "04-Nov-2025"
$helper.Copy()

foreach ($r in $rows) {
    $hCell = $ws.Cells.Item($r, 8)   # column H
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($r, 9)   # column I
    $iCell.PasteSpecial(-4163)       # xlPasteValues
}

# Remove the temporary helper row so it leaves no trace in the sheet.
$helper.EntireRow.Delete()
